$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value = 6

$ws.Range("F7").Value = -3
$ws.Range("F8").Value = -1
$ws.Range("F9").Value = -4
$ws.Range("F10").Value = 2
$ws.Range("F11").Value = 4
$ws.Range("F12").Value = -2
$ws.Range("F13").Value = -2
$ws.Range("F14").Value = -1
$ws.Range("F15").Value = 4
$ws.Range("F16").Value = -3
$ws.Range("F17").Value = -6
$ws.Range("F18").Value = 6
$ws.Range("F19").Value = -1
$ws.Range("F20").Value = 4
$ws.Range("F21").Value = 0
